$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value, derived from the recalculated NATMI
# ligand-receptor pair statistics (Flt3l-Flt3.xlsx), per "Natmi following
# Dr Hou advice" commit.
$updates = @{}
$updates["E2"] = 3
$updates["G2"] = 9.801662666666665
$updates["H2"] = 29.404988
$updates["I2"] = 0.3903913270717663
$updates["J2"] = 0.3903913270717663
$updates["K2"] = 3
$updates["M2"] = 0.6054233333333333
$updates["N2"] = 1.81627
$updates["O2"] = 0.3419053608852526
$updates["P2"] = 0.3419053608852526
$updates["Q2"] = 5.934155283862221
$updates["R2"] = 53.40739755476
$updates["S2"] = 0.133476887568945
$updates["T2"] = 0.1334768875689449
$updates["E3"] = 3
$updates["G3"] = 9.801662666666665
$updates["H3"] = 29.404988
$updates["I3"] = 0.3903913270717663
$updates["J3"] = 0.3903913270717663
$updates["K3"] = 3
$updates["M3"] = 1.165310333333333
$updates["N3"] = 3.495931
$updates["O3"] = 0.6580946391147473
$updates["P3"] = 0.6580946391147473
$updates["Q3"] = 11.42197878931422
$updates["R3"] = 102.797809103828
$updates["S3"] = 0.2569144395028213
$updates["T3"] = 0.2569144395028213
$updates["E4"] = 3
$updates["G4"] = 5.350150333333333
$updates["H4"] = 16.050451
$updates["I4"] = 0.2130916314602937
$updates["J4"] = 0.2130916314602937
$updates["K4"] = 3
$updates["M4"] = 0.6054233333333333
$updates["N4"] = 1.81627
$updates["O4"] = 0.3419053608852526
$updates["P4"] = 0.3419053608852526
$updates["Q4"] = 3.239105848641111
$updates["R4"] = 29.15195263777
$updates["S4"] = 0.07285717115605898
$updates["T4"] = 0.07285717115605897
$updates["E5"] = 3
$updates["G5"] = 5.350150333333333
$updates["H5"] = 16.050451
$updates["I5"] = 0.2130916314602937
$updates["J5"] = 0.2130916314602937
$updates["K5"] = 3
$updates["M5"] = 1.165310333333333
$updates["N5"] = 3.495931
$updates["O5"] = 0.6580946391147473
$updates["P5"] = 0.6580946391147473
$updates["Q5"] = 6.23458546832011
$updates["R5"] = 56.11126921488099
$updates["S5"] = 0.1402344603042347
$updates["T5"] = 0.1402344603042347
$updates["E6"] = 3
$updates["G6"] = 3.545239
$updates["H6"] = 10.635717
$updates["I6"] = 0.1412036513665554
$updates["J6"] = 0.1412036513665554
$updates["K6"] = 3
$updates["M6"] = 0.6054233333333333
$updates["N6"] = 1.81627
$updates["O6"] = 0.3419053608852526
$updates["P6"] = 0.3419053608852526
$updates["Q6"] = 2.146370412843333
$updates["R6"] = 19.31733371559
$updates["S6"] = 0.04827828537879753
$updates["T6"] = 0.04827828537879753
$updates["E7"] = 3
$updates["G7"] = 3.545239
$updates["H7"] = 10.635717
$updates["I7"] = 0.1412036513665554
$updates["J7"] = 0.1412036513665554
$updates["K7"] = 3
$updates["M7"] = 1.165310333333333
$updates["N7"] = 3.495931
$updates["O7"] = 0.6580946391147473
$updates["P7"] = 0.6580946391147473
$updates["Q7"] = 4.131303640836333
$updates["R7"] = 37.181732767527
$updates["S7"] = 0.0929253659877579
$updates["T7"] = 0.09292536598775787
$updates["E8"] = 3
$updates["G8"] = 5.087578
$updates["H8"] = 15.262734
$updates["I8"] = 0.2026336137597937
$updates["J8"] = 0.2026336137597937
$updates["K8"] = 3
$updates["M8"] = 0.6054233333333333
$updates["N8"] = 1.81627
$updates["O8"] = 0.3419053608852526
$updates["P8"] = 0.3419053608852526
$updates["Q8"] = 3.080138431353333
$updates["R8"] = 27.72124588218
$updates["S8"] = 0.06928151884002516
$updates["T8"] = 0.06928151884002515
$updates["E9"] = 3
$updates["G9"] = 5.087578
$updates["H9"] = 15.262734
$updates["I9"] = 0.2026336137597937
$updates["J9"] = 0.2026336137597937
$updates["K9"] = 3
$updates["M9"] = 1.165310333333333
$updates["N9"] = 3.495931
$updates["O9"] = 0.6580946391147473
$updates["P9"] = 0.6580946391147473
$updates["Q9"] = 5.928607215039333
$updates["R9"] = 53.35746493535399
$updates["S9"] = 0.1333520949197685
$updates["T9"] = 0.1333520949197685
$updates["E10"] = 3
$updates["G10"] = 1.322645666666667
$updates["H10"] = 3.967937
$updates["I10"] = 0.05267977634159087
$updates["J10"] = 0.05267977634159086
$updates["K10"] = 3
$updates["M10"] = 0.6054233333333333
$updates["N10"] = 1.81627
$updates["O10"] = 0.3419053608852526
$updates["P10"] = 0.3419053608852526
$updates["Q10"] = 0.8007605483322222
$updates["R10"] = 7.206844934989999
$updates["S10"] = 0.01801149794142602
$updates["T10"] = 0.01801149794142601
$updates["E11"] = 3
$updates["G11"] = 1.322645666666667
$updates["H11"] = 3.967937
$updates["I11"] = 0.05267977634159087
$updates["J11"] = 0.05267977634159086
$updates["K11"] = 3
$updates["M11"] = 1.165310333333333
$updates["N11"] = 3.495931
$updates["O11"] = 0.6580946391147473
$updates["P11"] = 0.6580946391147473
$updates["Q11"] = 1.541292662705222
$updates["R11"] = 13.871633964347
$updates["S11"] = 0.03466827840016485
$updates["T11"] = 0.03466827840016484

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
